# BA_tasks.xlsx - "feat: first work on eigenfaces"
#
# Summary of the edit:
#  - D23 text updated: "Google's universal Sentence Encoder: Alter shapes to
#    fix problem" -> "... (HOW?)"
#  - New row: D22 = "make TFIDF searchable"
#  - B19 text replaced with a longer, merged note about InferSent / Universal
#    Sentence Encoder / eigenfaces
#  - Row heights adjusted to fit the new/changed wrapped text
#  - Selection moved to B19 (single cell, not A19:B19) and the view scrolled
#    down a couple of rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content updates -------------------------------------------------

$ws.Range("D23").Value = "Google's universal Sentence Encoder: Alter shapes to fix problem (HOW?)"

$ws.Range("D22").Value = "make TFIDF searchable"

$ws.Range("B19").Value = "InferSent hypothese: project does not support training: https://github.com/facebookresearch/InferSent/issues/82 , Universal Sentence Encoder: the embedding uses n-grams of documents close to current doc (like a window) to embed it, cf. DAN in https://amitness.com/2020/06/universal-sentence-encoder/ -> unable to fix subtraction problem, started working on eigenfaces"

# --- row heights (auto-recalculated by Excel for wrapped text) -------

$ws.Rows.Item(19).RowHeight = 154
$ws.Rows.Item(22).RowHeight = 17
$ws.Rows.Item(23).RowHeight = 34

# --- view state: scroll + selection -----------------------------------

$win = $excel.ActiveWindow
$win.ScrollRow = 17
[void]$ws.Range("B19").Select()
